$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as exact text (preserve formatting/precision)
$textCells = @("D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D20","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data values
$ws.Range("D2").Value = "25.678.70"
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("D3").Value = "1.740.91"
$ws.Range("E3").Value = "  -5.52%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "238.91"
$ws.Range("E5").Value = "  -8.06%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "0.4909"
$ws.Range("E7").Value = "  -7.21%  "
$ws.Range("D8").Value = "41.90"
$ws.Range("E8").Value = "  -7.13%  "
$ws.Range("D9").Value = "0.2434"
$ws.Range("E9").Value = "  -22.78%  "
$ws.Range("D10").Value = "0.05997"
$ws.Range("E10").Value = "  -11.80%  "
$ws.Range("D11").Value = "1.734.11"
$ws.Range("E11").Value = "  -5.42%  "
$ws.Range("D12").Value = "0.06768"
$ws.Range("E12").Value = "  -12.74%  "
$ws.Range("D13").Value = "14.72"
$ws.Range("E13").Value = "  -21.41%  "
$ws.Range("D14").Value = "0.5902"
$ws.Range("E14").Value = "  -24.46%  "
$ws.Range("D15").Value = "4.412"
$ws.Range("E15").Value = "  -11.95%  "
$ws.Range("D16").Value = "76.63"
$ws.Range("E16").Value = "  -12.91%  "
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "25.707.92"
$ws.Range("E19").Value = "  -3.60%  "
$ws.Range("D20").Value = "11.51"
$ws.Range("E20").Value = "  -16.92%  "
$ws.Range("D21").Value = "0.000006355"
$ws.Range("E21").Value = "  -19.82%  "
$ws.Range("D22").Value = "1.953.76"
$ws.Range("E22").Value = "  -5.88%  "
$ws.Range("D23").Value = "3.939"
$ws.Range("E23").Value = "  -14.50%  "
$ws.Range("D24").Value = "5.136"
$ws.Range("E24").Value = "  -14.10%  "
$ws.Range("D25").Value = "7.843"
$ws.Range("E25").Value = "  -15.93%  "
$ws.Range("D26").Value = "135.31"
$ws.Range("E26").Value = "  -5.30%  "
$ws.Range("D27").Value = "1.843"
$ws.Range("E27").Value = "  -16.40%  "
$ws.Range("D28").Value = "1.452"
$ws.Range("E28").Value = "  -13.76%  "
$ws.Range("D29").Value = "14.46"
$ws.Range("E29").Value = "  -14.97%  "
$ws.Range("D30").Value = "100.61"
$ws.Range("E30").Value = "  -9.35%  "
$ws.Range("D31").Value = "0.08103"
$ws.Range("E31").Value = "  -7.14%  "
$ws.Range("D32").Value = "3.719"
$ws.Range("E32").Value = "  -11.17%  "
$ws.Range("D33").Value = "3.407"
$ws.Range("E33").Value = "  -16.42%  "
$ws.Range("D34").Value = "0.04343"
$ws.Range("E34").Value = "  -10.84%  "
$ws.Range("D35").Value = "0.9993"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "2.669"
$ws.Range("E36").Value = "  -6.61%  "
$ws.Range("D37").Value = "1.026"
$ws.Range("E37").Value = "  -10.14%  "
$ws.Range("E38").Value = "  -17.33%  "
$ws.Range("D39").Value = "2.737"
$ws.Range("E39").Value = "  -11.41%  "
$ws.Range("D40").Value = "2.045"
$ws.Range("E40").Value = "  -9.42%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "102.15"
$ws.Range("E42").Value = "  -7.16%  "
$ws.Range("D43").Value = "0.01489"
$ws.Range("E43").Value = "  -14.11%  "
$ws.Range("D44").Value = "0.7986"
$ws.Range("E44").Value = "  -10.72%  "
$ws.Range("D45").Value = "0.3821"
$ws.Range("E45").Value = "  -20.38%  "
$ws.Range("D46").Value = "5.099"
$ws.Range("E46").Value = "  -13.81%  "
$ws.Range("D47").Value = "6.044"
$ws.Range("E47").Value = "  -21.23%  "
$ws.Range("D48").Value = "0.05093"
$ws.Range("E48").Value = "  -12.43%  "
$ws.Range("D49").Value = "30.03"
$ws.Range("E49").Value = "  -13.62%  "
$ws.Range("D50").Value = "52.28"
$ws.Range("E50").Value = "  -12.54%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.1037"
$ws.Range("E51").Value = "  -16.39%  "
